# Update "想去人数" (F column) values on sheet "展览" (sheet1) and "全部类型" (sheet4)
# to reflect the latest generated output, per commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (first sheet) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 1615
$ws1.Range("F3").Value  = 9002
$ws1.Range("F6").Value  = 690
$ws1.Range("F7").Value  = 349
$ws1.Range("F8").Value  = 178
$ws1.Range("F9").Value  = 49
$ws1.Range("F10").Value = 80
$ws1.Range("F11").Value = 3841
$ws1.Range("F13").Value = 379
$ws1.Range("F15").Value = 4283
$ws1.Range("F18").Value = 1144
$ws1.Range("F21").Value = 5
$ws1.Range("F22").Value = 243
$ws1.Range("F24").Value = 2649
$ws1.Range("F25").Value = 114

# --- Sheet "全部类型" (fourth sheet) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 1615
$ws4.Range("F3").Value  = 9002
$ws4.Range("F6").Value  = 690
$ws4.Range("F7").Value  = 349
$ws4.Range("F8").Value  = 178
$ws4.Range("F9").Value  = 49
$ws4.Range("F10").Value = 80
$ws4.Range("F11").Value = 3841
$ws4.Range("F13").Value = 379
$ws4.Range("F15").Value = 4283
$ws4.Range("F18").Value = 1144
$ws4.Range("F21").Value = 5
$ws4.Range("F22").Value = 243
$ws4.Range("F24").Value = 2649
$ws4.Range("F26").Value = 114
